# Trade #4 closed at 2026-02-17 19:43:59 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1300.04
$wsSummary.Range("B4").Value = 0.04
$wsSummary.Range("B5").Value = 0.2
$wsSummary.Range("B6").Value = 4
$wsSummary.Range("B7").Value = 2
$wsSummary.Range("B9").Value = 50

# --- Strategy Status sheet ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.04
$wsStatus.Range("D4").Value = 4
$wsStatus.Range("E4").Value = 0.04
$wsStatus.Range("F4").Value = 0.04
$wsStatus.Range("G4").Value = 50

# --- Helper to append the newly closed trade as row 5 ---
function Add-TradeRow($ws) {
    $ws.Cells.Item(5, 1).Value = 4

    # "2026-02-17" looks like a date to Excel's auto type detection, which
    # would turn it into a date serial number + date number format. Force
    # the cell to Text first, assign the literal string, then restore the
    # cell style to Normal so no stray numFmt/style survives on the cell.
    $dateCell = $ws.Cells.Item(5, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $ws.Cells.Item(5, 3).Value = "19:43:53"
    $ws.Cells.Item(5, 4).Value = "MarketMaking"
    $ws.Cells.Item(5, 5).Value = "DOWN"
    $ws.Cells.Item(5, 6).Value = 0.7
    $ws.Cells.Item(5, 7).Value = 0.72
    $ws.Cells.Item(5, 8).Value = "CLOSED"
    $ws.Cells.Item(5, 9).Value = 2.8571
    $ws.Cells.Item(5, 10).Value = 0.02
    $ws.Cells.Item(5, 11).Value = 100.04
    $ws.Cells.Item(5, 12).Value = 0
    $ws.Cells.Item(5, 13).Value = 0
    $ws.Cells.Item(5, 14).Value = 0.6
    $ws.Cells.Item(5, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(5, 16).Value = "early_exit"
    $ws.Cells.Item(5, 17).Value = 0.13
}

# --- All Trades sheet ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

# --- MarketMaking sheet ---
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
